$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VO IDs range")

# Add the new "vaccine against pathogen/disease" intermediate term ID,
# replacing the reserved placeholder ID that had been entered in A13.
$ws.Range("A13").Value = "VO:0010458"

# Move the active selection to A14, matching where editing continued.
$ws.Range("A14").Select()
